# Weekly update: a new price record for "Cilantro" at "Macroferia Regional
# de Talca" is prepended to the data table. Insert a new row at row 55
# (pushing the existing rows 55-89 down to 56-90) and populate it with the
# new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("55:55").Insert()

$ws.Range("A55").Value = 5
$ws.Range("B55").Value = "Macroferia Regional de Talca"
$ws.Range("C55").Value = "Maule"
$ws.Range("D55").Value = 45062
$ws.Range("D55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E55").Value = 7
$ws.Range("F55").Value = 100112040
$ws.Range("G55").Value = "Cilantro"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 7000
$ws.Range("L55").Value = 7000
$ws.Range("M55").Value = 7000
$ws.Range("N55").Value = "`$/caja 36 atados"
$ws.Range("O55").Value = "Región del Maule"
$ws.Range("P55").Value = 194
$ws.Range("Q55").Value = 36
$ws.Range("R55").Value = "Hortaliza"
